$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing style of column D price cells, force Text format while
# writing so Excel does not auto-convert numeric-looking strings (e.g. "210.75")
# into real numbers, then restore the original style/format afterwards so the
# cell formatting is left exactly as it was.
$dStyle = $ws.Range("D2:D51").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.297.21'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.566.01'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '210.75'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '44.39'
$ws.Range("E8").Value = '  -4.21%  '
$ws.Range("D9").Value = '23.73'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").Value = '0.0585'
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = '0.0895'
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").Value = '1.790.10'
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '1.569.22'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '28.298.89'
$ws.Range("E16").Value = '  -0.87%  '
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '227.39'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").Value = '0.0₃0675'
$ws.Range("E21").Value = '  -2.60%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").Value = '8.90'
$ws.Range("E24").Value = '  -2.49%  '
$ws.Range("D25").Value = '2.04'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").Value = '150.41'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '14.87'
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").Value = '0.0478'
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").Value = '1.07'
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("D35").Value = '1.379.23'
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("E37").Value = '  -2.82%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").Value = '2.64'
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("D41").Value = '0.520'
$ws.Range("E41").Value = '  -2.74%  '
$ws.Range("D42").Value = '1.93'
$ws.Range("E42").Value = '  +3.55%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("D47").Value = '62.11'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '0.916'
$ws.Range("E48").Value = '  -6.32%  '
$ws.Range("D49").Value = '1.702.73'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '85.35'
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("E51").Value = '  -1.95%  '

# Restore original column D style/number format
$ws.Range("D2:D51").Style = $dStyle
